$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing number format for the Price (D) column, switch to Text so
# numeric-looking price strings (e.g. "580.11") are written as literal text,
# matching the workbook author's inline-string cells, then restore the format.
$priceRange = $ws.Range("D2:D51")
$origStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.621.70"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.249.52"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "580.11"
$ws.Range("D6").Value = "183.98"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "3.247.81"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "3.799.58"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D15").Value = "27.67"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "67.638.87"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "3.213.05"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "13.53"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "394.60"
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "71.47"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "9.58"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D31").Value = "5.54"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").Value = "22.63"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "7.01"
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "161.55"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D39").Value = "26.62"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "0.807"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("D44").Value = "0.0688"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "40.67"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "2.616.82"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "334.44"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -0.51%  "

# Restore the original (General) style/number format for the Price column.
$priceRange.Style = $origStyle

